$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.044.93'
$ws.Range("E2").Value = '  -0.02%  '
$ws.Range("D3").Value = '1.833.34'
$ws.Range("E3").Value = '  +0.17%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9990'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '244.12'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +1.36%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6344'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +1.91%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9998'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07551'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +0.56%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2948'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +1.05%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '22.94'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +0.85%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07734'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +1.41%  '
$ws.Range("D12").Value = '1.836.23'
$ws.Range("E12").Value = '  +0.40%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.007'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +1.04%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6713'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +1.16%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '83.22'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +1.24%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.000009698'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +6.60%  '
$ws.Range("E17").Value = '  +1.68%  '
$ws.Range("D18").Value = '29.084.96'
$ws.Range("E18").Value = '  +0.16%  '
$ws.Range("E19").Value = '  +2.17%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '226.33'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +0.72%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9990'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -0.11%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.201'
$ws.Range("D22").ClearFormats()
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.9995'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -0.08%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '160.53'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +0.68%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1403'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +3.04%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.556'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +1.74%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.93'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +0.61%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.498'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +0.05%  '
$ws.Range("E29").Value = '  +2.03%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.082'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +1.30%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.203'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -0.35%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.05387'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +3.24%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.866'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +1.85%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7457'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +1.76%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.144'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -0.81%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.656'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +0.49%  '
$ws.Range("D37").Value = '1.243.62'
$ws.Range("E37").Value = '  -2.27%  '
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01790'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +0.39%  '
$ws.Range("B39").Value = 'MXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.756'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +0.26%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.646'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +5.08%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9070'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +1.79%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9997'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -0.08%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '101.93'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +0.02%  '
$ws.Range("D44").Value = '1.984.58'
$ws.Range("E44").Value = '  +0.47%  '
$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '64.93'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +2.45%  '
$ws.Range("B46").Value = 'BabyDogeCoin'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00000000122'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +2.53%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5109'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -0.13%  '
$ws.Range("E48").Value = '  +3.31%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.074'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +2.34%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.779'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +2.09%  '
$ws.Range("B51").Value = 'RenderToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.649'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -0.97%  '
